# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Updates the "K" column (column G) values for rows 2-16 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 2
    4  = 1
    5  = 1
    6  = 1
    7  = 2
    8  = 1
    9  = 1
    10 = 0
    11 = 4
    12 = 2
    13 = 3
    14 = 0
    15 = 6
    16 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
